{"js": "// Apply the \"Dritte Woche\" weekly-report updates:\n//  1. Zweite Woche / \"Eigener Zeitaufwand\": extend the school-hours figure\n//     with a decimal value and an extra parenthetical school-period count.\n//  2. Dritte Woche / \"Eigene Leistung\": fill in the (until now empty) entry.\n//  3. Dritte Woche / \"Gruppenleistung\": fill in the (until now empty) entry.\n//  4. Dritte Woche / \"Eigener Zeitaufwand\": replace the old figure with the\n//     new home/school hour breakdown.\n\nconst body = context.document.body;\n\n// --- 1. Zweite Woche / Eigener Zeitaufwand ------------------------------\n// \" zuhause, 3 Stunden in der Schule\" -> \" zuhause, 1.5 Stunden in der\n// Schule (10,5 Schulstunden)\"\nconst zweiteWocheHit = body.search(\" zuhause, 3 Stunden in der Schule\", {\n  matchCase: true\n});\nzweiteWocheHit.load(\"text\");\nawait context.sync();\n\nif (zweiteWocheHit.items.length !== 1) {\n  throw new Error(\n    \"Expected exactly one match for the Zweite-Woche Zeitaufwand text, found \" +\n      zweiteWocheHit.items.length\n  );\n}\nzweiteWocheHit.items[0].insertText(\n  \" zuhause, 1.5 Stunden in der Schule (10,5 Schulstunden)\",\n  \"Replace\"\n);\nawait context.sync();\n\n// --- Locate the \"Dritte Woche\" section via its bold heading paragraph ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet drittWocheIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Dritte Woche\") === 0) {\n    drittWocheIndex = i;\n    break;\n  }\n}\nif (drittWocheIndex === -1) {\n  throw new Error('Could not locate the \"Dritte Woche\" heading paragraph.');\n}\n\nconst eigeneLeistungPara = paragraphs.items[drittWocheIndex + 1];\nconst gruppenleistungPara = paragraphs.items[drittWocheIndex + 2];\nconst zeitaufwandPara = paragraphs.items[drittWocheIndex + 3];\n\neigeneLeistungPara.load(\"text\");\ngruppenleistungPara.load(\"text\");\nzeitaufwandPara.load(\"text\");\nawait context.sync();\n\nif (eigeneLeistungPara.text !== \"Eigene Leistung:\") {\n  throw new Error(\n    'Expected \"Eigene Leistung:\" paragraph, found \"' +\n      eigeneLeistungPara.text +\n      '\"'\n  );\n}\nif (gruppenleistungPara.text !== \"Gruppenleistung:\") {\n  throw new Error(\n    'Expected \"Gruppenleistung:\" paragraph, found \"' +\n      gruppenleistungPara.text +\n      '\"'\n  );\n}\nif (zeitaufwandPara.text !== \"Eigener Zeitaufwand: 3 Stunden (Samstag)\") {\n  throw new Error(\n    'Expected \"Eigener Zeitaufwand: 3 Stunden (Samstag)\" paragraph, found \"' +\n      zeitaufwandPara.text +\n      '\"'\n  );\n}\n\n// --- 2. Dritte Woche / Eigene Leistung ----------------------------------\n// \"Eigene Leistung:\" -> \"Eigene Leistung: Multithreading, Mausinteraktion,\n// Erste Schritte f\u00fcr Gravitationssimulation\"\neigeneLeistungPara\n  .getRange(\"End\")\n  .insertText(\n    \" Multithreading, Mausinteraktion, Erste Schritte f\u00fcr Gravitationssimulation\",\n    \"End\"\n  );\nawait context.sync();\n\n// --- 3. Dritte Woche / Gruppenleistung ----------------------------------\n// \"Gruppenleistung:\" -> \"Gruppenleistung: Effizientere Darstellung der\n// Boids, realistischeres Verhalten, Klicken st\u00f6\u00dft Boids ab\"\ngruppenleistungPara\n  .getRange(\"End\")\n  .insertText(\n    \" Effizientere Darstellung der Boids, realistischeres Verhalten, Klicken st\u00f6\u00dft Boids ab\",\n    \"End\"\n  );\nawait context.sync();\n\n// --- 4. Dritte Woche / Eigener Zeitaufwand ------------------------------\n// \" 3 Stunden (Samstag)\" -> \" 4 Schulstunden zuhause, 2 Schulstunden in\n// der Schule\"\nconst zeitaufwandHit = zeitaufwandPara\n  .getRange()\n  .search(\" 3 Stunden (Samstag)\", { matchCase: true });\nzeitaufwandHit.load(\"text\");\nawait context.sync();\n\nif (zeitaufwandHit.items.length !== 1) {\n  throw new Error(\n    \"Expected exactly one match for the Dritte-Woche Zeitaufwand text, found \" +\n      zeitaufwandHit.items.length\n  );\n}\nzeitaufwandHit.items[0].insertText(\n  \" 4 Schulstunden zuhause, 2 Schulstunden in der Schule\",\n  \"Replace\"\n);\nawait context.sync();\n", "ps1": "# Apply the \"Dritte Woche\" weekly-report updates:\n#  1. Zweite Woche / \"Eigener Zeitaufwand\": extend the school-hours figure\n#     with a decimal value and an extra parenthetical school-period count.\n#  2. Dritte Woche / \"Eigene Leistung\": fill in the (until now empty) entry.\n#  3. Dritte Woche / \"Gruppenleistung\": fill in the (until now empty) entry.\n#  4. Dritte Woche / \"Eigener Zeitaufwand\": replace the old figure with the\n#     new home/school hour breakdown.\n\n$d = $word.ActiveDocument\n\n# --- 1. Zweite Woche / Eigener Zeitaufwand ------------------------------\n# \" zuhause, 3 Stunden in der Schule\" -> \" zuhause, 1.5 Stunden in der\n# Schule (10,5 Schulstunden)\"\n$zweiteWocheRange = $d.Content\n$found1 = $zweiteWocheRange.Find.Execute(\" zuhause, 3 Stunden in der Schule\")\nif (-not $found1) {\n    throw \"Could not find the Zweite-Woche Eigener-Zeitaufwand text.\"\n}\n$zweiteWocheRange.Text = \" zuhause, 1.5 Stunden in der Schule (10,5 Schulstunden)\"\n\n# --- Locate the \"Dritte Woche\" section via its bold heading paragraph ---\n$paragraphCount = $d.Paragraphs.Count\n$drittWocheIndex = -1\nfor ($i = 1; $i -le $paragraphCount; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text.StartsWith(\"Dritte Woche\")) {\n        $drittWocheIndex = $i\n        break\n    }\n}\nif ($drittWocheIndex -eq -1) {\n    throw 'Could not locate the \"Dritte Woche\" heading paragraph.'\n}\n\n$eigeneLeistungPara = $d.Paragraphs.Item($drittWocheIndex + 1)\n$gruppenleistungPara = $d.Paragraphs.Item($drittWocheIndex + 2)\n$zeitaufwandPara = $d.Paragraphs.Item($drittWocheIndex + 3)\n\nif ($eigeneLeistungPara.Range.Text.Trim() -ne \"Eigene Leistung:\") {\n    throw \"Expected 'Eigene Leistung:' paragraph, found '$($eigeneLeistungPara.Range.Text.Trim())'\"\n}\nif ($gruppenleistungPara.Range.Text.Trim() -ne \"Gruppenleistung:\") {\n    throw \"Expected 'Gruppenleistung:' paragraph, found '$($gruppenleistungPara.Range.Text.Trim())'\"\n}\nif ($zeitaufwandPara.Range.Text.Trim() -ne \"Eigener Zeitaufwand: 3 Stunden (Samstag)\") {\n    throw \"Expected 'Eigener Zeitaufwand: 3 Stunden (Samstag)' paragraph, found '$($zeitaufwandPara.Range.Text.Trim())'\"\n}\n\n# --- 2. Dritte Woche / Eigene Leistung ----------------------------------\n# \"Eigene Leistung:\" -> \"Eigene Leistung: Multithreading, Mausinteraktion,\n# Erste Schritte f\u00fcr Gravitationssimulation\"\n$eigeneLeistungRange = $eigeneLeistungPara.Range.Duplicate\n[void]$eigeneLeistungRange.MoveEnd(1, -1)\n[void]$eigeneLeistungRange.InsertAfter(\" Multithreading, Mausinteraktion, Erste Schritte f\u00fcr Gravitationssimulation\")\n\n# --- 3. Dritte Woche / Gruppenleistung ----------------------------------\n# \"Gruppenleistung:\" -> \"Gruppenleistung: Effizientere Darstellung der\n# Boids, realistischeres Verhalten, Klicken st\u00f6\u00dft Boids ab\"\n$gruppenleistungRange = $gruppenleistungPara.Range.Duplicate\n[void]$gruppenleistungRange.MoveEnd(1, -1)\n[void]$gruppenleistungRange.InsertAfter(\" Effizientere Darstellung der Boids, realistischeres Verhalten, Klicken st\u00f6\u00dft Boids ab\")\n\n# --- 4. Dritte Woche / Eigener Zeitaufwand -------------------------------\n# \" 3 Stunden (Samstag)\" -> \" 4 Schulstunden zuhause, 2 Schulstunden in\n# der Schule\"\n$zeitaufwandRange = $zeitaufwandPara.Range.Duplicate\n$found4 = $zeitaufwandRange.Find.Execute(\" 3 Stunden (Samstag)\")\nif (-not $found4) {\n    throw \"Could not find the Dritte-Woche Eigener-Zeitaufwand text.\"\n}\n$zeitaufwandRange.Text = \" 4 Schulstunden zuhause, 2 Schulstunden in der Schule\"\n"}
